# Re-fill previously empty Genus (C) / species-epithet (D) cells with the
# literal text "NA" so downstream TEXTJOIN results read "..._NA_NA" instead
# of "..._NA" for rows that had no species/epithet recorded.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cellsToFill = @(
    "C4","C10","C15","C18","C20","C31","C50","C68","C74","C80",
    "C107","C118","C119","C123","C124","C152","C157","C159","C213","C214",
    "C217","C368","D2","D4","D6","D7","D9","D10","D12","D13",
    "D14","D15","D18","D19","D20","D21","D23","D24","D25","D27",
    "D28","D29","D31","D32","D33","D34","D36","D37","D38","D45",
    "D47","D48","D50","D54","D56","D57","D58","D61","D62","D64",
    "D66","D68","D69","D71","D74","D78","D80","D81","D85","D86",
    "D87","D88","D90","D91","D92","D95","D103","D107","D109","D110",
    "D113","D116","D118","D119","D123","D124","D125","D130","D133","D134",
    "D137","D138","D139","D140","D148","D151","D152","D155","D157","D158",
    "D159","D166","D169","D170","D171","D181","D182","D183","D189","D191",
    "D194","D195","D196","D197","D199","D206","D207","D208","D213","D214",
    "D217","D219","D225","D228","D230","D232","D240","D245","D246","D247",
    "D248","D254","D269","D270","D295","D300","D303","D305","D336","D368"
)

foreach ($addr in $cellsToFill) {
    $ws.Range($addr).Value = "NA"
}

# Restore the sheet's scroll/selection state: frozen pane top-left at A2,
# active cell F4 (matches the saved view after the edit).
$ws.Range("F4").Select() | Out-Null
